# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp.
# - Refresh COVID-19 per-country metrics (Casos totales/Nuevos casos/Casos
#   activos/Recuperados/Casos criticos/Muertes hoy/Muertes) for the rows
#   whose figures moved between the two snapshots.
# - A handful of neighbouring countries swapped rank (their updated figures
#   pushed them past/behind a neighbour in the sorted list), so those rows
#   also get their country name (column A) corrected to match the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1)
$ws.Range('A1').Value = 'Datos actualizados a 24 de Agosto de 2020 a las 17:19'

# Each entry: row number, new country name ($null = unchanged), then the
# new B..H values (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes).
$rowUpdates = @(
    @(4,   $null,               5880879, 6733,  3167319, 2532860, 0, 96,  180700),
    @(6,   $null,               3126078, 20893, 2355823, 712364,  0, 199, 57891),
    @(13,  $null,               399568,  1903,  372464,  16188,   0, 64,  10916),
    @(15,  $null,               342154,  0,     256789,  78286,   0, 94,  7079),
    @(23,  $null,               235367,  878,   209600,  16433,   0, 2,   9334),
    @(27,  $null,               125001,  105,   111190,  4737,    0, 1,   9074),
    @(37,  $null,               86721,   0,     0,       0,       0, 6,   5813),
    @(41,  $null,               80960,   432,   72925,   7517,    0, 3,   518),
    @(62,  $null,               39156,   210,   35283,   3596,    0, 4,   277),
    @(64,  $null,               35426,   152,   33104,   1803,    0, 1,   519),
    @(65,  'Moldavia',          33828,   350,   23570,   9313,    0, 5,   945),
    @(66,  'Costa Rica',        33820,   0,     10518,   22947,   0, 0,   355),
    @(68,  $null,               32557,   193,   18895,   13108,   0, 6,   554),
    @(69,  $null,               30714,   57,    29028,   985,     0, 3,   701),
    @(91,  $null,               9842,    12,    8953,    638,     0, 0,   251),
    @(98,  'Tayikistan',        8346,    35,    7142,    1137,    0, 1,   67),
    @(99,  'Croacia',           8311,    136,   5926,    2212,    0, 2,   173),
    @(108, $null,               5383,    1,     5273,    50,      0, 0,   60),
    @(126, $null,               2959,    6,     2811,    136,     0, 0,   12),
    @(133, 'Uganda',            2362,    99,    1248,    1092,    0, 2,   22),
    @(134, 'Estonia',           2275,    3,     2025,    186,     0, 1,   64),
    @(149, $null,               1451,    30,    878,     553,     0, 0,   20),
    @(150, $null,               1421,    10,    1137,    266,     0, 1,   18),
    @(154, $null,               1290,    4,     819,     389,     0, 0,   82),
    @(159, 'Trinidad yTobago',  1031,    24,    165,     851,     0, 1,   15),
    @(160, 'Vietnam',           1022,    6,     587,     408,     0, 0,   27),
    @(161, 'Lesoto',            1015,    0,     472,     513,     0, 0,   30),
    @(171, 'Birmania',          474,     24,    341,     127,     0, 0,   6),
    @(172, 'Martinica',         464,     0,     98,      350,     0, 0,   16)
)

$cols = @('B', 'C', 'D', 'E', 'F', 'G', 'H')

foreach ($u in $rowUpdates) {
    $r = $u[0]
    $name = $u[1]
    if ($null -ne $name) {
        $ws.Range("A$r").Value = $name
    }
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $u[2 + $i]
    }
}
